$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.091.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.11%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.926.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.11%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.15%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'325.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.93%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.16%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4603"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.80%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3825"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.65%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07771"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.73%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9816"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +2.07%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'22.76"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +3.37%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.935.75"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +2.31%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'5.696"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.89%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'6.976"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.62%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.07044"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.62%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +0.20%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'84.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.60%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.000009541"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.65%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'16.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.10%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'1.003"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.14%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'29.095.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.30%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.345"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.89%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +1.08%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.077"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.05%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'157.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.23%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.74%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'5.679"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.68%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'118.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.12%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.837"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.54%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.09350"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.99%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.8607"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.92%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'5.119"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.05%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'1.244"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.75%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'3.016"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.26%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.161"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.53%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.05696"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.56%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'3.208"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +18.94%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'1.004"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.18%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.02050"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +1.23%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'7.501"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.64%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.5516"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.61%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1757"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.48%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +2.67%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'2.190"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +6.47%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.000002748"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -7.47%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.5193"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.85%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'11.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.28%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.06918"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.71%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -0.74%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'1.770"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.50%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +0.22%  "
$ws.Range("E51").Style = "Normal"
